$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 50 values (per MV data refresh) ---
$ws.Cells.Item(50, 7).Value = 1542    # G50
$ws.Cells.Item(50, 8).Value = 330     # H50
$ws.Cells.Item(50, 32).Value = 3093   # AF50
$ws.Cells.Item(50, 34).Value = 961    # AH50

# --- Append new row 51: quarter 01-04-2021 ---
# Write the period label as a genuine text string (not an auto-converted
# date serial) by building it via a formula and pasting the *value*,
# mirroring how the rest of column A stores its "dd-mm-yyyy"-look text.
$ws.Cells.Item(53, 1).Formula = "=""01-04-2021"""
$ws.Cells.Item(53, 1).Copy()
$ws.Cells.Item(51, 1).PasteSpecial(-4163)
$ws.Cells.Item(53, 1).Clear()

$ws.Cells.Item(51, 2).Value = 0
$ws.Cells.Item(51, 3).Value = 6
$ws.Cells.Item(51, 4).Value = -2
$ws.Cells.Item(51, 5).Value = 4
$ws.Cells.Item(51, 6).Value = 9350
$ws.Cells.Item(51, 7).Value = 570
$ws.Cells.Item(51, 8).Value = 279
$ws.Cells.Item(51, 9).Value = -1594
$ws.Cells.Item(51, 10).Value = 444
$ws.Cells.Item(51, 11).Value = -2462
$ws.Cells.Item(51, 12).Value = 400
$ws.Cells.Item(51, 13).Value = 57
$ws.Cells.Item(51, 14).Value = -596
$ws.Cells.Item(51, 15).Value = 64
$ws.Cells.Item(51, 16).Value = 0
$ws.Cells.Item(51, 17).Value = 31
$ws.Cells.Item(51, 18).Value = -38
$ws.Cells.Item(51, 19).Value = 36
$ws.Cells.Item(51, 20).Value = 269
$ws.Cells.Item(51, 21).Value = 170
$ws.Cells.Item(51, 22).Value = -604
$ws.Cells.Item(51, 23).Value = 210
$ws.Cells.Item(51, 24).Value = 2515
$ws.Cells.Item(51, 25).Value = 367
$ws.Cells.Item(51, 26).Value = -1270
$ws.Cells.Item(51, 27).Value = 245
$ws.Cells.Item(51, 28).Value = 0
$ws.Cells.Item(51, 29).Value = 1
$ws.Cells.Item(51, 30).Value = -1
$ws.Cells.Item(51, 31).Value = 2
$ws.Cells.Item(51, 32).Value = 3754
$ws.Cells.Item(51, 33).Value = 6888
$ws.Cells.Item(51, 34).Value = 911
$ws.Cells.Item(51, 35).Value = -4105
$ws.Cells.Item(51, 36).Value = 1005
